$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new observation is prepended before the current row 434,
# pushing the existing rows 434-453 down to 435-454.
$ws.Rows(434).Insert()

$ws.Range("A434").Value = 5
$ws.Range("B434").Value = "Macroferia Regional de Talca"
$ws.Range("C434").Value = "Maule"
$ws.Range("D434").Value = 44939
$ws.Range("E434").Value = 7
$ws.Range("F434").Value = 100114013
$ws.Range("G434").Value = "Zanahoria"
$ws.Range("H434").Value = "Sin especificar"
$ws.Range("I434").Value = "Primera"
$ws.Range("J434").Value = 400
$ws.Range("K434").Value = 10000
$ws.Range("L434").Value = 10000
$ws.Range("M434").Value = 10000
$ws.Range("N434").Value = '$/saco 20 kilos'
$ws.Range("O434").Value = "Región de Ñuble"
$ws.Range("P434").Value = 500
$ws.Range("Q434").Value = 20
$ws.Range("R434").Value = "Hortaliza"
